# Updates the cryptocurrency price/volume table (and for rows 33-34, swaps the
# Monero / EthereumClassic entries) to reflect the latest scrape, per the
# "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: values are written with a leading apostrophe so Excel stores them as
# plain text (matching the workbook's original inline-string cells) instead of
# auto-converting numeric-looking strings (e.g. "0.999", "7.03") into numbers.

# Row 2: Bitcoin
$ws.Range("D2").Value = "'58.601.16"
$ws.Range("E2").Value = "'  -1.61%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "'2.628.74"
$ws.Range("E3").Value = "'  +0.72%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.13%  "

# Row 5: BNB
$ws.Range("D5").Value = "'535.68"
$ws.Range("E5").Value = "'  -0.22%  "

# Row 6: Solana
$ws.Range("D6").Value = "'142.86"
$ws.Range("E6").Value = "'  +0.72%  "

# Row 7: USDC
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "'  -0.08%  "

# Row 8: XRP
$ws.Range("E8").Value = "'  -0.09%  "

# Row 9: LidoStakedEther
$ws.Range("D9").Value = "'2.636.40"
$ws.Range("E9").Value = "'  +0.73%  "

# Row 10: Toncoin
$ws.Range("D10").Value = "'7.03"
$ws.Range("E10").Value = "'  +8.35%  "

# Row 11: Dogecoin
$ws.Range("D11").Value = "'0.101"
$ws.Range("E11").Value = "'  -1.74%  "

# Row 12: Cardano
$ws.Range("E12").Value = "'  -0.23%  "

# Row 13: TRON
$ws.Range("E13").Value = "'  +0.99%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "'3.098.93"
$ws.Range("E14").Value = "'  +0.97%  "

# Row 15: WrappedBTC
$ws.Range("D15").Value = "'58.549.22"
$ws.Range("E15").Value = "'  -1.50%  "

# Row 16: Avalanche
$ws.Range("D16").Value = "'20.88"
$ws.Range("E16").Value = "'  +1.00%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "'2.637.62"
$ws.Range("E17").Value = "'  +0.95%  "

# Row 18: ShibaInu
$ws.Range("E18").Value = "'  -0.98%  "

# Row 19: Polkadot
$ws.Range("D19").Value = "'4.41"
$ws.Range("E19").Value = "'  +1.03%  "

# Row 20: BitcoinCash
$ws.Range("D20").Value = "'334.56"
$ws.Range("E20").Value = "'  -2.06%  "

# Row 21: Chainlink
$ws.Range("E21").Value = "'  +0.48%  "

# Row 22: Uniswap
$ws.Range("D22").Value = "'6.24"
$ws.Range("E22").Value = "'  -2.21%  "

# Row 23: Dai
$ws.Range("E23").Value = "'  -0.01%  "

# Row 24: Litecoin
$ws.Range("D24").Value = "'66.40"
$ws.Range("E24").Value = "'  -1.47%  "

# Row 25: Polygon
$ws.Range("E25").Value = "'  +1.43%  "

# Row 26: Kaspa
$ws.Range("E26").Value = "'  -0.84%  "

# Row 27: Binance-PegBSC-USD
$ws.Range("D27").Value = "'0.997"
$ws.Range("E27").Value = "'  -0.20%  "

# Row 28: InternetComputer(DFINITY)
$ws.Range("D28").Value = "'7.14"
$ws.Range("E28").Value = "'  -1.40%  "

# Row 29: PEPE
$ws.Range("D29").Value = "'0.0₃0739"
$ws.Range("E29").Value = "'  -1.34%  "

# Row 31: PancakeSwap
$ws.Range("D31").Value = "'1.65"
$ws.Range("E31").Value = "'  -1.47%  "

# Row 32: Aptos
$ws.Range("E32").Value = "'  +0.00%  "

# Row 33: Monero
$ws.Range("B33").Value = "'EthereumClassic"
$ws.Range("C33").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'18.74"
$ws.Range("E33").Value = "'  -0.61%  "

# Row 34: EthereumClassic
$ws.Range("B34").Value = "'Monero"
$ws.Range("C34").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").Value = "'150.60"
$ws.Range("E34").Value = "'  +0.68%  "

# Row 35: NEARProtocol
$ws.Range("E35").Value = "'  -2.13%  "

# Row 36: OKB
$ws.Range("D36").Value = "'37.17"
$ws.Range("E36").Value = "'  -0.17%  "

# Row 37: ImmutableX
$ws.Range("D37").Value = "'1.11"
$ws.Range("E37").Value = "'  -0.89%  "

# Row 38: SuiNetwork
$ws.Range("D38").Value = "'0.843"
$ws.Range("E38").Value = "'  +0.44%  "

# Row 39: Stacks
$ws.Range("E39").Value = "'  -3.47%  "

# Row 40: Fetch.AI
$ws.Range("D40").Value = "'0.814"
$ws.Range("E40").Value = "'  -2.09%  "

# Row 41: Filecoin
$ws.Range("E41").Value = "'  +1.02%  "

# Row 42: Bittensor
$ws.Range("D42").Value = "'281.76"
$ws.Range("E42").Value = "'  +2.45%  "

# Row 43: FirstDigitalUSD
$ws.Range("E43").Value = "'  -0.05%  "

# Row 44: Mantle
$ws.Range("E44").Value = "'  +0.14%  "

# Row 45: WhiteBITCoin
$ws.Range("E45").Value = "'  -0.32%  "

# Row 46: EnergySwap
$ws.Range("D46").Value = "'19.11"
$ws.Range("E46").Value = "'  +2.93%  "

# Row 47: Hedera
$ws.Range("E47").Value = "'  +1.44%  "

# Row 48: Stellar
$ws.Range("D48").Value = "'0.0936"
$ws.Range("E48").Value = "'  -2.03%  "

# Row 49: VeChain
$ws.Range("E49").Value = "'  +0.63%  "

# Row 50: Maker
$ws.Range("D50").Value = "'1.946.84"
$ws.Range("E50").Value = "'  -0.16%  "

# Row 51: RenderToken
$ws.Range("D51").Value = "'4.46"
$ws.Range("E51").Value = "'  -1.49%  "
